$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = Fecha (date serial numbers), Column J = Volumen
# Update per-row values as described by the diff (weekly logic re-shuffle)

$ws.Range("D2").Value2 = 44630
$ws.Range("J2").Value2 = 60

$ws.Range("D3").Value2 = 44659

$ws.Range("D4").Value2 = 44649
$ws.Range("J4").Value2 = 60

$ws.Range("D5").Value2 = 44645

$ws.Range("D6").Value2 = 44651
$ws.Range("J6").Value2 = 60

$ws.Range("D7").Value2 = 44642

$ws.Range("D8").Value2 = 44628
$ws.Range("J8").Value2 = 60

$ws.Range("D10").Value2 = 44635

$ws.Range("D11").Value2 = 44637
$ws.Range("J11").Value2 = 100

$ws.Range("D12").Value2 = 44658
$ws.Range("J12").Value2 = 80

$ws.Range("D13").Value2 = 44656
$ws.Range("J13").Value2 = 100

$ws.Range("D14").Value2 = 44664
$ws.Range("J14").Value2 = 160
